# Natmi following Dr Hou advice
#
# The underlying NATMI ligand-receptor analysis was re-run with a third
# "Sending cluster" / "Target cluster" category (ECs) added alongside the
# existing FAPs and sCs clusters. This grows the result table from
# 2 senders x 3 targets = 6 data rows to 3 senders x 3 targets = 9 data
# rows (rows 2-10, below the header in row 1), and every numeric column is
# recomputed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one full data row: Sending cluster, Ligand symbol,
# Receptor symbol, Target cluster, then the 16 numeric measure columns
# (E..T). The leading "," on each row forces PowerShell to keep this as
# an array-of-arrays instead of flattening everything into one list.
$rows = @(
  ,@("ECs",  "Lamc2", "Itga6", "ECs",  2, 0.6666666666666666, 0.4421816666666667, 1.326545,   0.06026482003168283, 0.06026482003168283, 3, 1, 103.4766596666667, 310.429979, 0.877785331764719,  0.8777853317647188, 45.75548183250611,  411.799336492555,   0.05289957504525179,  0.05289957504525179)
  ,@("ECs",  "Lamc2", "Itga6", "FAPs", 2, 0.6666666666666666, 0.4421816666666667, 1.326545,   0.06026482003168283, 0.06026482003168283, 3, 1, 0.8265796666666668, 2.479739,   0.007011818020336602, 0.0070118180203366,  0.3654983746394445,  3.289485371755001,  0.0004225659510904959, 0.0004225659510904958)
  ,@("ECs",  "Lamc2", "Itga6", "sCs",  2, 0.6666666666666666, 0.4421816666666667, 1.326545,   0.06026482003168283, 0.06026482003168283, 3, 1, 13.58054833333333, 40.741645,  0.1152028502149446,  0.1152028502149446, 6.005069496280556,  54.045625466525,    0.00694267903534055,  0.006942679035340548)
  ,@("FAPs", "Lamc2", "Itga6", "ECs",  3, 1,                  6.369908666666666,  19.109726,  0.8681531333236113,  0.8681531333236113,  3, 1, 103.4766596666667, 310.429979, 0.877785331764719,  0.8777853317647188, 659.1368712084171,  5932.231840875754,  0.7620520861570464,   0.7620520861570463)
  ,@("FAPs", "Lamc2", "Itga6", "FAPs", 3, 1,                  6.369908666666666,  19.109726,  0.8681531333236113,  0.8681531333236113,  3, 1, 0.8265796666666668, 2.479739,   0.007011818020336602, 0.0070118180203366,  5.265236982390445,  47.387132841514,    0.006087331784650182, 0.00608733178465018)
  ,@("FAPs", "Lamc2", "Itga6", "sCs",  3, 1,                  6.369908666666666,  19.109726,  0.8681531333236113,  0.8681531333236113,  3, 1, 13.58054833333333, 40.741645,  0.1152028502149446,  0.1152028502149446, 86.50685252658555,  778.5616727392699,  0.1000137153819148,   0.1000137153819148)
  ,@("sCs",  "Lamc2", "Itga6", "ECs",  3, 1,                  0.5252196666666666, 1.575659,   0.07158204664470585, 0.07158204664470584, 3, 1, 103.4766596666667, 310.429979, 0.877785331764719,  0.8777853317647188, 54.34797669790677,  489.131790281161,   0.06283367056242072,  0.0628336705624207)
  ,@("sCs",  "Lamc2", "Itga6", "FAPs", 3, 1,                  0.5252196666666666, 1.575659,   0.07158204664470585, 0.07158204664470584, 3, 1, 0.8265796666666668, 2.479739,   0.007011818020336602, 0.0070118180203366,  0.4341358970001112,  3.907223073001,     0.0005019202845959236, 0.0005019202845959234)
  ,@("sCs",  "Lamc2", "Itga6", "sCs",  3, 1,                  0.5252196666666666, 1.575659,   0.07158204664470585, 0.07158204664470584, 3, 1, 13.58054833333333, 40.741645,  0.1152028502149446,  0.1152028502149446, 7.132771068783888,  64.19493961905499,  0.008246455797689227, 0.008246455797689223)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowValues = $rows[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
